$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 - NroSiniestro
$ws.Range("E3").Value = "'1220194200684"

# Row 2 - NroSiniestro: trailing space added (leading "'" forces text, matching
# the existing quotePrefix style already applied to these cells)
$ws.Range("E2").Value = "'1120194100442 "

# Row 5 - NroSiniestro changed
$ws.Range("E5").Value = "'0420172008629    "

# Row 6 - NroSiniestro changed (trailing spaces 3 -> 4)
$ws.Range("E6").Value = "'1220170301442    "

# Row 7 - NroSiniestro changed (trailing spaces 3 -> 4)
$ws.Range("E7").Value = "'1120170200942    "

# Row 5 - Usuario changed
$ws.Range("C5").Value = "mpimpignano"

# Update the active selection to E5
$ws.Range("E5").Select()
